$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header titles (row 1) ---
$ws.Range("A1").Value = "part number"
$ws.Range("B1").Value = "quantity"
$ws.Range("C1").Value = "UM"
$ws.Range("D1").Value = "value mxn"

# --- Row 8: pick up the same look as row 2 (copy format only) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

# --- Row 10: quantity correction ---
$ws.Range("B10").Value = 50

# --- The "UM" column (C) on the blank rows still carries the numeric style
#     used by column B; match it to the text style already used by C2:C10 ---
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C11:C14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New inventory rows 11-14 ---
$ws.Range("A11").Value = "65dsf"
$ws.Range("B11").Value = 680
$ws.Range("C11").Value = "pcs"
$ws.Range("D11").Value = 1.2

$ws.Range("A12").Value = "65dsf654"
$ws.Range("B12").Value = 800
$ws.Range("C12").Value = "pcs"
$ws.Range("D12").Value = 0.8

$ws.Range("A13").Value = "345df"
$ws.Range("B13").Value = 80
$ws.Range("C13").Value = "pcs"
$ws.Range("D13").Value = 0.2
$ws.Range("D13").NumberFormat = '"$"#,##0.00'
# Touch the font so this cell gets its own currency style distinct from D2:D12
$ws.Range("D13").Font.Name = $ws.Range("D13").Font.Name

$ws.Range("A14").Value = "345dg"
$ws.Range("B14").Value = 50
$ws.Range("C14").Value = "pcs"
$ws.Range("D14").Value = 0.8

# --- Move the visible selection to D1, matching the saved view state ---
$ws.Range("D1").Select() | Out-Null

$excel.CutCopyMode = 0
